# Apply the edits described by the target diff:
#  - Add a new worksheet "Kmeans-2" at the end of the workbook, with
#    l/h performance column data (A/B) and a small TP/TN/FP/FN summary
#    table (E/F).
#  - Add a TP/TN/FP/FN summary table (S/T) to sheet "2" (the second tab).
#  - Make sheet "2" the active tab / sheet, update selections on the
#    first two sheets.
#  - Nudge sheet "2"'s page setup (orientation) so a <pageSetup> element
#    is emitted.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet "2" (2nd tab): add the TP:/TN:/FP:/FN: counts in columns S:T
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("S2").Value = "TP:"
$ws2.Range("T2").Value = 18
$ws2.Range("T2").NumberFormat = "0"

$ws2.Range("S3").Value = "TN:"
$ws2.Range("T3").Value = 12
$ws2.Range("T3").NumberFormat = "0"

$ws2.Range("S4").Value = "FP:"
$ws2.Range("T4").Value = 17
$ws2.Range("T4").NumberFormat = "0"

$ws2.Range("S5").Value = "FN:"
$ws2.Range("T5").Value = 20
$ws2.Range("T5").NumberFormat = "0"

# Give sheet 2 a page setup (shows up as <pageSetup orientation="portrait".../>)
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 2. Add the new "Kmeans-2" worksheet at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Kmeans-2"

$newSheet.Range("A1").Value = "Performance"
$newSheet.Range("B1").Value = "Cluster"

$aVals = @("l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","l","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h","h")
$bVals = @(0,1,1,1,1,1,0,1,1,1,1,1,0,0,1,0,0,1,0,0,1,1,0,1,0,1,0,1,0,0,0,0,0,1,1,0,1,0,1,1,0,1,1,1,0,1,0,1,0,0,1,1,1,1,1,1,0,1,1,1,0,1,1,1,0,0,1)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $aVals[$i]
    $newSheet.Cells.Item($row, 2).Value = $bVals[$i]
}

# TP:/TN:/FP:/FN: summary table in columns E:F
$newSheet.Range("E2").Value = "TP:"
$newSheet.Range("F2").Value = 19
$newSheet.Range("E3").Value = "TN:"
$newSheet.Range("F3").Value = 11
$newSheet.Range("E4").Value = "FP:"
$newSheet.Range("F4").Value = 16
$newSheet.Range("E5").Value = "FN:"
$newSheet.Range("F5").Value = 21

$newSheet.Range("E2:E5").Select()

# ---------------------------------------------------------------------
# 3. Selections / active sheet bookkeeping
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:A67").Select()

$ws2.Activate()
$ws2.Range("V18").Select()
